# Remove the "noon_close" column (G) from the sheet, shifting the
# Fr_afternoon/Sa_morning/Sa_afternoon/Su_morning/Su_afternoon columns
# one position to the left (G..K instead of H..L).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1:G2").EntireColumn.Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)

# Update the selection to match the new layout.
$ws.Range("J2:K2").Select()

$wb.Save()
